$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7202794431374137
$ws.Range("C2").Value = 0.9788435059368326
$ws.Range("D2").Value = 0.6477765369869576
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 AdaBoostRegressor(learning_rate=0.5, n_estimators=100))])"
$ws.Range("G2").Value = 0.1242467469831657
$ws.Range("H2").Value = 0.991
